$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("TODOS")
$ws.Cells.Item(9,1).Value = "19:53"
$ws.Cells.Item(9,2).Value = "16_SANTA ANA"
$ws.Cells.Item(9,3).Value = 0
$ws.Cells.Item(9,4).Value = "🚌"
$ws.Cells.Item(10,1).Value = "18:31"
$ws.Cells.Item(10,2).Value = "15_ABASTO"
$ws.Cells.Item(10,3).Value = 1
$ws.Cells.Item(10,4).Value = "🚌"
$ws.Cells.Item(11,1).Value = "19:37"
$ws.Cells.Item(11,2).Value = "15_ABASTO"
$ws.Cells.Item(11,3).Value = 1
$ws.Cells.Item(11,4).Value = "🚌"
$ws.Cells.Item(12,1).Value = "18:35"
$ws.Cells.Item(12,2).Value = "23_HERNANDEZ"
$ws.Cells.Item(12,3).Value = 5
$ws.Cells.Item(12,4).Value = "🚌"
$ws.Cells.Item(13,1).Value = "19:41"
$ws.Cells.Item(13,2).Value = "14_ABASTO"
$ws.Cells.Item(13,3).Value = 5
$ws.Cells.Item(13,4).Value = "🚌"
$ws.Cells.Item(14,1).Value = "18:40"
$ws.Cells.Item(14,2).Value = "14_ABASTO"
$ws.Cells.Item(14,3).Value = 10
$ws.Cells.Item(14,4).Value = "📅"
$ws.Cells.Item(15,1).Value = "18:40"
$ws.Cells.Item(15,2).Value = "15_ABASTO"
$ws.Cells.Item(15,3).Value = 10
$ws.Cells.Item(15,4).Value = "🚌"
$ws.Cells.Item(16,1).Value = "19:51"
$ws.Cells.Item(16,2).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(16,3).Value = 15
$ws.Cells.Item(16,4).Value = "🚌"
$ws.Cells.Item(17,1).Value = "19:00"
$ws.Cells.Item(17,2).Value = "16_SANTA ANA"
$ws.Cells.Item(17,3).Value = 17
$ws.Cells.Item(17,4).Value = "🚌"
$ws.Cells.Item(18,1).Value = "18:48"
$ws.Cells.Item(18,2).Value = "14X44_ABASTO"
$ws.Cells.Item(18,3).Value = 18
$ws.Cells.Item(18,4).Value = "🚌"
$ws.Cells.Item(19,1).Value = "20:12"
$ws.Cells.Item(19,2).Value = "23_HERNANDEZ"
$ws.Cells.Item(19,3).Value = 19
$ws.Cells.Item(19,4).Value = "🚌"
$ws.Cells.Item(20,1).Value = "19:04"
$ws.Cells.Item(20,2).Value = "23_HERNANDEZ"
$ws.Cells.Item(20,3).Value = 21
$ws.Cells.Item(20,4).Value = "🚌"
$ws.Cells.Item(21,1).Value = "18:52"
$ws.Cells.Item(21,2).Value = "215A_LA PLATA"
$ws.Cells.Item(21,3).Value = 22
$ws.Cells.Item(21,4).Value = "🚌"
$ws.Cells.Item(22,1).Value = "20:01"
$ws.Cells.Item(22,2).Value = "14_ABASTO"
$ws.Cells.Item(22,3).Value = 25
$ws.Cells.Item(22,4).Value = "📅"
$ws.Cells.Item(23,1).Value = "18:56"
$ws.Cells.Item(23,2).Value = "10_OLMOS"
$ws.Cells.Item(23,3).Value = 26
$ws.Cells.Item(23,4).Value = "🚌"
$ws.Cells.Item(24,1).Value = "19:10"
$ws.Cells.Item(24,2).Value = "14_ABASTO"
$ws.Cells.Item(24,3).Value = 27
$ws.Cells.Item(24,4).Value = "🚌"
$ws.Cells.Item(25,1).Value = "18:58"
$ws.Cells.Item(25,2).Value = "215A_EL PATO"
$ws.Cells.Item(25,3).Value = 28
$ws.Cells.Item(25,4).Value = "📅"
$ws.Cells.Item(26,1).Value = "19:12"
$ws.Cells.Item(26,2).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(26,3).Value = 29
$ws.Cells.Item(26,4).Value = "🚌"
$ws.Cells.Item(27,1).Value = "20:22"
$ws.Cells.Item(27,2).Value = "16_SANTA ANA"
$ws.Cells.Item(27,3).Value = 29
$ws.Cells.Item(27,4).Value = "🚌"
$ws.Cells.Item(28,1).Value = "19:30"
$ws.Cells.Item(28,2).Value = "16_SANTA ANA"
$ws.Cells.Item(28,3).Value = 30
$ws.Cells.Item(28,4).Value = "🚌"
$ws.Cells.Item(29,1).Value = "19:01"
$ws.Cells.Item(29,2).Value = "16_SANTA ANA"
$ws.Cells.Item(29,3).Value = 31
$ws.Cells.Item(29,4).Value = "🚌"
$ws.Cells.Item(30,1).Value = "19:16"
$ws.Cells.Item(30,2).Value = "15_ABASTO"
$ws.Cells.Item(30,3).Value = 33
$ws.Cells.Item(30,4).Value = "🚌"
$ws.Cells.Item(31,1).Value = "19:16"
$ws.Cells.Item(31,2).Value = "27_EL RETIRO"
$ws.Cells.Item(31,3).Value = 33
$ws.Cells.Item(31,4).Value = "📅"
$ws.Cells.Item(32,1).Value = "19:04"
$ws.Cells.Item(32,2).Value = "11_ETCHEVERRY"
$ws.Cells.Item(32,3).Value = 34
$ws.Cells.Item(32,4).Value = "🚌"
$ws.Cells.Item(33,1).Value = "19:05"
$ws.Cells.Item(33,2).Value = "23_HERNANDEZ"
$ws.Cells.Item(33,3).Value = 35
$ws.Cells.Item(33,4).Value = "🚌"
$ws.Cells.Item(34,1).Value = "20:11"
$ws.Cells.Item(34,2).Value = "10_OLMOS"
$ws.Cells.Item(34,3).Value = 35
$ws.Cells.Item(34,4).Value = "🚌"
$ws.Cells.Item(35,1).Value = "20:11"
$ws.Cells.Item(35,2).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(35,3).Value = 35
$ws.Cells.Item(35,4).Value = "📅"
$ws.Cells.Item(36,1).Value = "19:20"
$ws.Cells.Item(36,2).Value = "14_ABASTO"
$ws.Cells.Item(36,3).Value = 37
$ws.Cells.Item(36,4).Value = "📅"
$ws.Cells.Item(37,1).Value = "19:20"
$ws.Cells.Item(37,2).Value = "16_SANTA ANA"
$ws.Cells.Item(37,3).Value = 37
$ws.Cells.Item(37,4).Value = "🚌"
$ws.Cells.Item(38,1).Value = "20:13"
$ws.Cells.Item(38,2).Value = "23_HERNANDEZ"
$ws.Cells.Item(38,3).Value = 37
$ws.Cells.Item(38,4).Value = "🚌"
$ws.Cells.Item(39,1).Value = "19:10"
$ws.Cells.Item(39,2).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(39,3).Value = 40
$ws.Cells.Item(39,4).Value = "🚌"
$ws.Cells.Item(40,1).Value = "19:10"
$ws.Cells.Item(40,2).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(40,3).Value = 40
$ws.Cells.Item(40,4).Value = "🚌"
$ws.Cells.Item(41,1).Value = "19:12"
$ws.Cells.Item(41,2).Value = "10_OLMOS"
$ws.Cells.Item(41,3).Value = 42
$ws.Cells.Item(41,4).Value = "🚌"
$ws.Cells.Item(42,1).Value = "19:28"
$ws.Cells.Item(42,2).Value = "15_ABASTO"
$ws.Cells.Item(42,3).Value = 45
$ws.Cells.Item(42,4).Value = "🚌"
$ws.Cells.Item(43,1).Value = "19:17"
$ws.Cells.Item(43,2).Value = "27_EL RETIRO"
$ws.Cells.Item(43,3).Value = 47
$ws.Cells.Item(43,4).Value = "🚌"
$ws.Cells.Item(44,1).Value = "20:24"
$ws.Cells.Item(44,2).Value = "215A_EL PATO"
$ws.Cells.Item(44,3).Value = 48
$ws.Cells.Item(44,4).Value = "📅"
$ws.Cells.Item(45,1).Value = "19:21"
$ws.Cells.Item(45,2).Value = "16_SANTA ANA"
$ws.Cells.Item(45,3).Value = 51
$ws.Cells.Item(45,4).Value = "🚌"
$ws.Cells.Item(46,1).Value = "19:21"
$ws.Cells.Item(46,2).Value = "26_HERNANDEZ"
$ws.Cells.Item(46,3).Value = 51
$ws.Cells.Item(46,4).Value = "🚌"
$ws.Cells.Item(47,1).Value = "19:34"
$ws.Cells.Item(47,2).Value = "23_HERNANDEZ"
$ws.Cells.Item(47,3).Value = 51
$ws.Cells.Item(47,4).Value = "🚌"
$ws.Cells.Item(48,1).Value = "20:52"
$ws.Cells.Item(48,2).Value = "23_HERNANDEZ"
$ws.Cells.Item(48,3).Value = 59
$ws.Cells.Item(48,4).Value = "🚌"
$ws.Cells.Item(49,1).Value = "19:30"
$ws.Cells.Item(49,2).Value = "225_GOMEZ"
$ws.Cells.Item(49,3).Value = 60
$ws.Cells.Item(49,4).Value = "📅"
$ws.Cells.Item(50,1).Value = "20:44"
$ws.Cells.Item(50,2).Value = "11_ETCHEVERRY"
$ws.Cells.Item(50,3).Value = 68
$ws.Cells.Item(50,4).Value = "🚌"
$ws.Cells.Item(51,1).Value = "20:09"
$ws.Cells.Item(51,2).Value = "15_ABASTO"
$ws.Cells.Item(51,3).Value = 69
$ws.Cells.Item(51,4).Value = "🚌"
$ws.Cells.Item(52,1).Value = "19:40"
$ws.Cells.Item(52,2).Value = "14_ABASTO"
$ws.Cells.Item(52,3).Value = 70
$ws.Cells.Item(52,4).Value = "🚌"
$ws.Cells.Item(53,1).Value = "19:40"
$ws.Cells.Item(53,2).Value = "215C_EL PATO"
$ws.Cells.Item(53,3).Value = 70
$ws.Cells.Item(53,4).Value = "🚌"
$ws.Cells.Item(54,1).Value = "20:10"
$ws.Cells.Item(54,2).Value = "10_OLMOS"
$ws.Cells.Item(54,3).Value = 70
$ws.Cells.Item(54,4).Value = "🚌"
$ws.Cells.Item(55,1).Value = "19:50"
$ws.Cells.Item(55,2).Value = "11X44_ETCHEVERRY"
$ws.Cells.Item(55,3).Value = 80
$ws.Cells.Item(55,4).Value = "🚌"
$ws.Cells.Item(56,1).Value = "19:50"
$ws.Cells.Item(56,2).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(56,3).Value = 80
$ws.Cells.Item(56,4).Value = "🚌"
$ws.Cells.Item(57,1).Value = "20:56"
$ws.Cells.Item(57,2).Value = "10_OLMOS"
$ws.Cells.Item(57,3).Value = 80
$ws.Cells.Item(57,4).Value = "🚌"
$ws.Cells.Item(58,1).Value = "19:51"
$ws.Cells.Item(58,2).Value = "81_EL PELIGRO"
$ws.Cells.Item(58,3).Value = 81
$ws.Cells.Item(58,4).Value = "🚌"
$ws.Cells.Item(59,1).Value = "20:57"
$ws.Cells.Item(59,2).Value = "27_EL RETIRO"
$ws.Cells.Item(59,3).Value = 81
$ws.Cells.Item(59,4).Value = "🚌"
$ws.Cells.Item(60,1).Value = "19:54"
$ws.Cells.Item(60,2).Value = "215C_LA PLATA"
$ws.Cells.Item(60,3).Value = 84
$ws.Cells.Item(60,4).Value = "🚌"
$ws.Cells.Item(61,1).Value = "21:04"
$ws.Cells.Item(61,2).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(61,3).Value = 88
$ws.Cells.Item(61,4).Value = "🚌"
$ws.Cells.Item(62,1).Value = "19:59"
$ws.Cells.Item(62,2).Value = "17_ROMERO"
$ws.Cells.Item(62,3).Value = 89
$ws.Cells.Item(62,4).Value = "📅"
$ws.Cells.Item(63,1).Value = "21:23"
$ws.Cells.Item(63,2).Value = "10_OLMOS"
$ws.Cells.Item(63,3).Value = 90
$ws.Cells.Item(63,4).Value = "🚌"
$ws.Cells.Item(64,1).Value = "21:23"
$ws.Cells.Item(64,2).Value = "15_ABASTO"
$ws.Cells.Item(64,3).Value = 90
$ws.Cells.Item(64,4).Value = "🚌"
$ws.Cells.Item(65,1).Value = "21:08"
$ws.Cells.Item(65,2).Value = "215B_EL PATO"
$ws.Cells.Item(65,3).Value = 92
$ws.Cells.Item(65,4).Value = "🚌"
$ws.Cells.Item(66,1).Value = "20:10"
$ws.Cells.Item(66,2).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(66,3).Value = 100
$ws.Cells.Item(66,4).Value = "🚌"
$ws.Cells.Item(67,1).Value = "21:37"
$ws.Cells.Item(67,2).Value = "14_ABASTO"
$ws.Cells.Item(67,3).Value = 104
$ws.Cells.Item(67,4).Value = "🚌"
$ws.Cells.Item(68,1).Value = "21:21"
$ws.Cells.Item(68,2).Value = "26_HERNANDEZ"
$ws.Cells.Item(68,3).Value = 105
$ws.Cells.Item(68,4).Value = "📅"
$ws.Cells.Item(69,1).Value = "21:38"
$ws.Cells.Item(69,2).Value = "17_ROMERO"
$ws.Cells.Item(69,3).Value = 105
$ws.Cells.Item(69,4).Value = "📅"
$ws.Cells.Item(70,1).Value = "20:31"
$ws.Cells.Item(70,2).Value = "225_GOMEZ"
$ws.Cells.Item(70,3).Value = 108
$ws.Cells.Item(70,4).Value = "📅"
$ws.Cells.Item(71,1).Value = "20:21"
$ws.Cells.Item(71,2).Value = "26_HERNANDEZ"
$ws.Cells.Item(71,3).Value = 111
$ws.Cells.Item(71,4).Value = "🚌"
$ws.Cells.Item(72,1).Value = "20:22"
$ws.Cells.Item(72,2).Value = "11_ETCHEVERRY"
$ws.Cells.Item(72,3).Value = 112
$ws.Cells.Item(72,4).Value = "🚌"
$ws.Cells.Item(73,1).Value = "20:23"
$ws.Cells.Item(73,2).Value = "215A_EL PATO"
$ws.Cells.Item(73,3).Value = 113
$ws.Cells.Item(73,4).Value = "🚌"
$ws.Cells.Item(74,1).Value = "21:29"
$ws.Cells.Item(74,2).Value = "215C_LA PLATA"
$ws.Cells.Item(74,3).Value = 113
$ws.Cells.Item(74,4).Value = "🚌"
$ws.Cells.Item(75,1).Value = "21:47"
$ws.Cells.Item(75,2).Value = "215A_EL PATO"
$ws.Cells.Item(75,3).Value = 114
$ws.Cells.Item(75,4).Value = "📅"
$ws.Cells.Item(76,1).Value = "20:39"
$ws.Cells.Item(76,2).Value = "215A_LA PLATA"
$ws.Cells.Item(76,3).Value = 116
$ws.Cells.Item(76,4).Value = "📅"
$ws.Cells.Item(77,1).Value = "20:56"
$ws.Cells.Item(77,2).Value = "27_EL RETIRO"
$ws.Cells.Item(77,3).Value = 116
$ws.Cells.Item(77,4).Value = "🚌"

$ws215 = $wb.Worksheets.Item("215")
$ws215.Cells.Item(2,1).Value = "18:52"
$ws215.Cells.Item(2,2).Value = "215A_LA PLATA"
$ws215.Cells.Item(2,3).Value = 22
$ws215.Cells.Item(2,4).Value = "🚌"
$ws215.Cells.Item(3,1).Value = "18:58"
$ws215.Cells.Item(3,2).Value = "215A_EL PATO"
$ws215.Cells.Item(3,3).Value = 28
$ws215.Cells.Item(3,4).Value = "📅"
$ws215.Cells.Item(4,1).Value = "19:12"
$ws215.Cells.Item(4,2).Value = "215B_LP-P MOR-1 Y 57"
$ws215.Cells.Item(4,3).Value = 29
$ws215.Cells.Item(4,4).Value = "🚌"
$ws215.Cells.Item(5,1).Value = "19:10"
$ws215.Cells.Item(5,2).Value = "215B_LP-P MOR-1 Y 57"
$ws215.Cells.Item(5,3).Value = 40
$ws215.Cells.Item(5,4).Value = "🚌"
$ws215.Cells.Item(6,1).Value = "20:24"
$ws215.Cells.Item(6,2).Value = "215A_EL PATO"
$ws215.Cells.Item(6,3).Value = 48
$ws215.Cells.Item(6,4).Value = "📅"
$ws215.Cells.Item(7,1).Value = "19:40"
$ws215.Cells.Item(7,2).Value = "215C_EL PATO"
$ws215.Cells.Item(7,3).Value = 70
$ws215.Cells.Item(7,4).Value = "🚌"
$ws215.Cells.Item(8,1).Value = "19:54"
$ws215.Cells.Item(8,2).Value = "215C_LA PLATA"
$ws215.Cells.Item(8,3).Value = 84
$ws215.Cells.Item(8,4).Value = "🚌"
$ws215.Cells.Item(9,1).Value = "21:08"
$ws215.Cells.Item(9,2).Value = "215B_EL PATO"
$ws215.Cells.Item(9,3).Value = 92
$ws215.Cells.Item(9,4).Value = "🚌"
$ws215.Cells.Item(10,1).Value = "20:23"
$ws215.Cells.Item(10,2).Value = "215A_EL PATO"
$ws215.Cells.Item(10,3).Value = 113
$ws215.Cells.Item(10,4).Value = "🚌"
$ws215.Cells.Item(11,1).Value = "21:29"
$ws215.Cells.Item(11,2).Value = "215C_LA PLATA"
$ws215.Cells.Item(11,3).Value = 113
$ws215.Cells.Item(11,4).Value = "🚌"
$ws215.Cells.Item(12,1).Value = "21:47"
$ws215.Cells.Item(12,2).Value = "215A_EL PATO"
$ws215.Cells.Item(12,3).Value = 114
$ws215.Cells.Item(12,4).Value = "📅"
$ws215.Cells.Item(13,1).Value = "20:39"
$ws215.Cells.Item(13,2).Value = "215A_LA PLATA"
$ws215.Cells.Item(13,3).Value = 116
$ws215.Cells.Item(13,4).Value = "📅"

$wsc = $wb.Worksheets.Item("COMBINADAS")
$wsc.Cells.Item(9,1).Value = "19:53"
$wsc.Cells.Item(9,2).Value = "16_SANTA ANA"
$wsc.Cells.Item(9,3).Value = 0
$wsc.Cells.Item(9,4).Value = "🚌"
$wsc.Cells.Item(10,1).Value = "18:31"
$wsc.Cells.Item(10,2).Value = "15_ABASTO"
$wsc.Cells.Item(10,3).Value = 1
$wsc.Cells.Item(10,4).Value = "🚌"
$wsc.Cells.Item(11,1).Value = "19:37"
$wsc.Cells.Item(11,2).Value = "15_ABASTO"
$wsc.Cells.Item(11,3).Value = 1
$wsc.Cells.Item(11,4).Value = "🚌"
$wsc.Cells.Item(12,1).Value = "18:35"
$wsc.Cells.Item(12,2).Value = "23_HERNANDEZ"
$wsc.Cells.Item(12,3).Value = 5
$wsc.Cells.Item(12,4).Value = "🚌"
$wsc.Cells.Item(13,1).Value = "19:41"
$wsc.Cells.Item(13,2).Value = "14_ABASTO"
$wsc.Cells.Item(13,3).Value = 5
$wsc.Cells.Item(13,4).Value = "🚌"
$wsc.Cells.Item(14,1).Value = "18:40"
$wsc.Cells.Item(14,2).Value = "14_ABASTO"
$wsc.Cells.Item(14,3).Value = 10
$wsc.Cells.Item(14,4).Value = "📅"
$wsc.Cells.Item(15,1).Value = "18:40"
$wsc.Cells.Item(15,2).Value = "15_ABASTO"
$wsc.Cells.Item(15,3).Value = 10
$wsc.Cells.Item(15,4).Value = "🚌"
$wsc.Cells.Item(16,1).Value = "19:51"
$wsc.Cells.Item(16,2).Value = "16_P MOR-SANTA ANA"
$wsc.Cells.Item(16,3).Value = 15
$wsc.Cells.Item(16,4).Value = "🚌"
$wsc.Cells.Item(17,1).Value = "19:00"
$wsc.Cells.Item(17,2).Value = "16_SANTA ANA"
$wsc.Cells.Item(17,3).Value = 17
$wsc.Cells.Item(17,4).Value = "🚌"
$wsc.Cells.Item(18,1).Value = "18:48"
$wsc.Cells.Item(18,2).Value = "14X44_ABASTO"
$wsc.Cells.Item(18,3).Value = 18
$wsc.Cells.Item(18,4).Value = "🚌"
$wsc.Cells.Item(19,1).Value = "20:12"
$wsc.Cells.Item(19,2).Value = "23_HERNANDEZ"
$wsc.Cells.Item(19,3).Value = 19
$wsc.Cells.Item(19,4).Value = "🚌"
$wsc.Cells.Item(20,1).Value = "19:04"
$wsc.Cells.Item(20,2).Value = "23_HERNANDEZ"
$wsc.Cells.Item(20,3).Value = 21
$wsc.Cells.Item(20,4).Value = "🚌"
$wsc.Cells.Item(21,1).Value = "18:52"
$wsc.Cells.Item(21,2).Value = "215A_LA PLATA"
$wsc.Cells.Item(21,3).Value = 22
$wsc.Cells.Item(21,4).Value = "🚌"
$wsc.Cells.Item(22,1).Value = "20:01"
$wsc.Cells.Item(22,2).Value = "14_ABASTO"
$wsc.Cells.Item(22,3).Value = 25
$wsc.Cells.Item(22,4).Value = "📅"
$wsc.Cells.Item(23,1).Value = "18:56"
$wsc.Cells.Item(23,2).Value = "10_OLMOS"
$wsc.Cells.Item(23,3).Value = 26
$wsc.Cells.Item(23,4).Value = "🚌"
$wsc.Cells.Item(24,1).Value = "19:10"
$wsc.Cells.Item(24,2).Value = "14_ABASTO"
$wsc.Cells.Item(24,3).Value = 27
$wsc.Cells.Item(24,4).Value = "🚌"
$wsc.Cells.Item(25,1).Value = "18:58"
$wsc.Cells.Item(25,2).Value = "215A_EL PATO"
$wsc.Cells.Item(25,3).Value = 28
$wsc.Cells.Item(25,4).Value = "📅"
$wsc.Cells.Item(26,1).Value = "19:12"
$wsc.Cells.Item(26,2).Value = "215B_LP-P MOR-1 Y 57"
$wsc.Cells.Item(26,3).Value = 29
$wsc.Cells.Item(26,4).Value = "🚌"
$wsc.Cells.Item(27,1).Value = "20:22"
$wsc.Cells.Item(27,2).Value = "16_SANTA ANA"
$wsc.Cells.Item(27,3).Value = 29
$wsc.Cells.Item(27,4).Value = "🚌"
$wsc.Cells.Item(28,1).Value = "19:30"
$wsc.Cells.Item(28,2).Value = "16_SANTA ANA"
$wsc.Cells.Item(28,3).Value = 30
$wsc.Cells.Item(28,4).Value = "🚌"
$wsc.Cells.Item(29,1).Value = "19:01"
$wsc.Cells.Item(29,2).Value = "16_SANTA ANA"
$wsc.Cells.Item(29,3).Value = 31
$wsc.Cells.Item(29,4).Value = "🚌"
$wsc.Cells.Item(30,1).Value = "19:16"
$wsc.Cells.Item(30,2).Value = "15_ABASTO"
$wsc.Cells.Item(30,3).Value = 33
$wsc.Cells.Item(30,4).Value = "🚌"
$wsc.Cells.Item(31,1).Value = "19:16"
$wsc.Cells.Item(31,2).Value = "27_EL RETIRO"
$wsc.Cells.Item(31,3).Value = 33
$wsc.Cells.Item(31,4).Value = "📅"
$wsc.Cells.Item(32,1).Value = "19:04"
$wsc.Cells.Item(32,2).Value = "11_ETCHEVERRY"
$wsc.Cells.Item(32,3).Value = 34
$wsc.Cells.Item(32,4).Value = "🚌"
$wsc.Cells.Item(33,1).Value = "19:05"
$wsc.Cells.Item(33,2).Value = "23_HERNANDEZ"
$wsc.Cells.Item(33,3).Value = 35
$wsc.Cells.Item(33,4).Value = "🚌"
$wsc.Cells.Item(34,1).Value = "20:11"
$wsc.Cells.Item(34,2).Value = "10_OLMOS"
$wsc.Cells.Item(34,3).Value = 35
$wsc.Cells.Item(34,4).Value = "🚌"
$wsc.Cells.Item(35,1).Value = "20:11"
$wsc.Cells.Item(35,2).Value = "16_P MOR-167 Y 521"
$wsc.Cells.Item(35,3).Value = 35
$wsc.Cells.Item(35,4).Value = "📅"
$wsc.Cells.Item(36,1).Value = "19:20"
$wsc.Cells.Item(36,2).Value = "14_ABASTO"
$wsc.Cells.Item(36,3).Value = 37
$wsc.Cells.Item(36,4).Value = "📅"
$wsc.Cells.Item(37,1).Value = "19:20"
$wsc.Cells.Item(37,2).Value = "16_SANTA ANA"
$wsc.Cells.Item(37,3).Value = 37
$wsc.Cells.Item(37,4).Value = "🚌"
$wsc.Cells.Item(38,1).Value = "20:13"
$wsc.Cells.Item(38,2).Value = "23_HERNANDEZ"
$wsc.Cells.Item(38,3).Value = 37
$wsc.Cells.Item(38,4).Value = "🚌"
$wsc.Cells.Item(39,1).Value = "19:10"
$wsc.Cells.Item(39,2).Value = "16_P MOR-SANTA ANA"
$wsc.Cells.Item(39,3).Value = 40
$wsc.Cells.Item(39,4).Value = "🚌"
$wsc.Cells.Item(40,1).Value = "19:10"
$wsc.Cells.Item(40,2).Value = "215B_LP-P MOR-1 Y 57"
$wsc.Cells.Item(40,3).Value = 40
$wsc.Cells.Item(40,4).Value = "🚌"
$wsc.Cells.Item(41,1).Value = "19:12"
$wsc.Cells.Item(41,2).Value = "10_OLMOS"
$wsc.Cells.Item(41,3).Value = 42
$wsc.Cells.Item(41,4).Value = "🚌"
$wsc.Cells.Item(42,1).Value = "19:28"
$wsc.Cells.Item(42,2).Value = "15_ABASTO"
$wsc.Cells.Item(42,3).Value = 45
$wsc.Cells.Item(42,4).Value = "🚌"
$wsc.Cells.Item(43,1).Value = "19:17"
$wsc.Cells.Item(43,2).Value = "27_EL RETIRO"
$wsc.Cells.Item(43,3).Value = 47
$wsc.Cells.Item(43,4).Value = "🚌"
$wsc.Cells.Item(44,1).Value = "20:24"
$wsc.Cells.Item(44,2).Value = "215A_EL PATO"
$wsc.Cells.Item(44,3).Value = 48
$wsc.Cells.Item(44,4).Value = "📅"
$wsc.Cells.Item(45,1).Value = "19:21"
$wsc.Cells.Item(45,2).Value = "16_SANTA ANA"
$wsc.Cells.Item(45,3).Value = 51
$wsc.Cells.Item(45,4).Value = "🚌"
$wsc.Cells.Item(46,1).Value = "19:21"
$wsc.Cells.Item(46,2).Value = "26_HERNANDEZ"
$wsc.Cells.Item(46,3).Value = 51
$wsc.Cells.Item(46,4).Value = "🚌"
$wsc.Cells.Item(47,1).Value = "19:34"
$wsc.Cells.Item(47,2).Value = "23_HERNANDEZ"
$wsc.Cells.Item(47,3).Value = 51
$wsc.Cells.Item(47,4).Value = "🚌"
$wsc.Cells.Item(48,1).Value = "20:52"
$wsc.Cells.Item(48,2).Value = "23_HERNANDEZ"
$wsc.Cells.Item(48,3).Value = 59
$wsc.Cells.Item(48,4).Value = "🚌"
$wsc.Cells.Item(49,1).Value = "19:30"
$wsc.Cells.Item(49,2).Value = "225_GOMEZ"
$wsc.Cells.Item(49,3).Value = 60
$wsc.Cells.Item(49,4).Value = "📅"
$wsc.Cells.Item(50,1).Value = "20:44"
$wsc.Cells.Item(50,2).Value = "11_ETCHEVERRY"
$wsc.Cells.Item(50,3).Value = 68
$wsc.Cells.Item(50,4).Value = "🚌"
$wsc.Cells.Item(51,1).Value = "20:09"
$wsc.Cells.Item(51,2).Value = "15_ABASTO"
$wsc.Cells.Item(51,3).Value = 69
$wsc.Cells.Item(51,4).Value = "🚌"
$wsc.Cells.Item(52,1).Value = "19:40"
$wsc.Cells.Item(52,2).Value = "14_ABASTO"
$wsc.Cells.Item(52,3).Value = 70
$wsc.Cells.Item(52,4).Value = "🚌"
$wsc.Cells.Item(53,1).Value = "19:40"
$wsc.Cells.Item(53,2).Value = "215C_EL PATO"
$wsc.Cells.Item(53,3).Value = 70
$wsc.Cells.Item(53,4).Value = "🚌"
$wsc.Cells.Item(54,1).Value = "20:10"
$wsc.Cells.Item(54,2).Value = "10_OLMOS"
$wsc.Cells.Item(54,3).Value = 70
$wsc.Cells.Item(54,4).Value = "🚌"
$wsc.Cells.Item(55,1).Value = "19:50"
$wsc.Cells.Item(55,2).Value = "11X44_ETCHEVERRY"
$wsc.Cells.Item(55,3).Value = 80
$wsc.Cells.Item(55,4).Value = "🚌"
$wsc.Cells.Item(56,1).Value = "19:50"
$wsc.Cells.Item(56,2).Value = "16_P MOR-SANTA ANA"
$wsc.Cells.Item(56,3).Value = 80
$wsc.Cells.Item(56,4).Value = "🚌"
$wsc.Cells.Item(57,1).Value = "20:56"
$wsc.Cells.Item(57,2).Value = "10_OLMOS"
$wsc.Cells.Item(57,3).Value = 80
$wsc.Cells.Item(57,4).Value = "🚌"
$wsc.Cells.Item(58,1).Value = "19:51"
$wsc.Cells.Item(58,2).Value = "81_EL PELIGRO"
$wsc.Cells.Item(58,3).Value = 81
$wsc.Cells.Item(58,4).Value = "🚌"
$wsc.Cells.Item(59,1).Value = "20:57"
$wsc.Cells.Item(59,2).Value = "27_EL RETIRO"
$wsc.Cells.Item(59,3).Value = 81
$wsc.Cells.Item(59,4).Value = "🚌"
$wsc.Cells.Item(60,1).Value = "19:54"
$wsc.Cells.Item(60,2).Value = "215C_LA PLATA"
$wsc.Cells.Item(60,3).Value = 84
$wsc.Cells.Item(60,4).Value = "🚌"
$wsc.Cells.Item(61,1).Value = "21:04"
$wsc.Cells.Item(61,2).Value = "84_COLONIA URQUIZA-ESC 49"
$wsc.Cells.Item(61,3).Value = 88
$wsc.Cells.Item(61,4).Value = "🚌"
$wsc.Cells.Item(62,1).Value = "19:59"
$wsc.Cells.Item(62,2).Value = "17_ROMERO"
$wsc.Cells.Item(62,3).Value = 89
$wsc.Cells.Item(62,4).Value = "📅"
$wsc.Cells.Item(63,1).Value = "21:23"
$wsc.Cells.Item(63,2).Value = "10_OLMOS"
$wsc.Cells.Item(63,3).Value = 90
$wsc.Cells.Item(63,4).Value = "🚌"
$wsc.Cells.Item(64,1).Value = "21:23"
$wsc.Cells.Item(64,2).Value = "15_ABASTO"
$wsc.Cells.Item(64,3).Value = 90
$wsc.Cells.Item(64,4).Value = "🚌"
$wsc.Cells.Item(65,1).Value = "21:08"
$wsc.Cells.Item(65,2).Value = "215B_EL PATO"
$wsc.Cells.Item(65,3).Value = 92
$wsc.Cells.Item(65,4).Value = "🚌"
$wsc.Cells.Item(66,1).Value = "20:10"
$wsc.Cells.Item(66,2).Value = "16_P MOR-167 Y 521"
$wsc.Cells.Item(66,3).Value = 100
$wsc.Cells.Item(66,4).Value = "🚌"
$wsc.Cells.Item(67,1).Value = "21:37"
$wsc.Cells.Item(67,2).Value = "14_ABASTO"
$wsc.Cells.Item(67,3).Value = 104
$wsc.Cells.Item(67,4).Value = "🚌"
$wsc.Cells.Item(68,1).Value = "21:21"
$wsc.Cells.Item(68,2).Value = "26_HERNANDEZ"
$wsc.Cells.Item(68,3).Value = 105
$wsc.Cells.Item(68,4).Value = "📅"
$wsc.Cells.Item(69,1).Value = "21:38"
$wsc.Cells.Item(69,2).Value = "17_ROMERO"
$wsc.Cells.Item(69,3).Value = 105
$wsc.Cells.Item(69,4).Value = "📅"
$wsc.Cells.Item(70,1).Value = "20:31"
$wsc.Cells.Item(70,2).Value = "225_GOMEZ"
$wsc.Cells.Item(70,3).Value = 108
$wsc.Cells.Item(70,4).Value = "📅"
$wsc.Cells.Item(71,1).Value = "20:21"
$wsc.Cells.Item(71,2).Value = "26_HERNANDEZ"
$wsc.Cells.Item(71,3).Value = 111
$wsc.Cells.Item(71,4).Value = "🚌"
$wsc.Cells.Item(72,1).Value = "20:22"
$wsc.Cells.Item(72,2).Value = "11_ETCHEVERRY"
$wsc.Cells.Item(72,3).Value = 112
$wsc.Cells.Item(72,4).Value = "🚌"
$wsc.Cells.Item(73,1).Value = "20:23"
$wsc.Cells.Item(73,2).Value = "215A_EL PATO"
$wsc.Cells.Item(73,3).Value = 113
$wsc.Cells.Item(73,4).Value = "🚌"
$wsc.Cells.Item(74,1).Value = "21:29"
$wsc.Cells.Item(74,2).Value = "215C_LA PLATA"
$wsc.Cells.Item(74,3).Value = 113
$wsc.Cells.Item(74,4).Value = "🚌"
$wsc.Cells.Item(75,1).Value = "21:47"
$wsc.Cells.Item(75,2).Value = "215A_EL PATO"
$wsc.Cells.Item(75,3).Value = 114
$wsc.Cells.Item(75,4).Value = "📅"
$wsc.Cells.Item(76,1).Value = "20:39"
$wsc.Cells.Item(76,2).Value = "215A_LA PLATA"
$wsc.Cells.Item(76,3).Value = 116
$wsc.Cells.Item(76,4).Value = "📅"
$wsc.Cells.Item(77,1).Value = "20:56"
$wsc.Cells.Item(77,2).Value = "27_EL RETIRO"
$wsc.Cells.Item(77,3).Value = 116
$wsc.Cells.Item(77,4).Value = "🚌"

Write-Output "done"